$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.868.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.80%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07358"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8682"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.758.81"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.353"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07085"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.483"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008690"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.88%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.903.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.329"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.79%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.036.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.900"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.163"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.303"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.62"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08926"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7660"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.157"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.497"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.904"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05273"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.948"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.262"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5300"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.349"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.407"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4920"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.30%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.50%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.665"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06274"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.01%  "
